# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    3  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    4  = @(3.230985683306322, 1.667794583268128, 337.1190423067083,  8.660232485948974, 350.6780550592317)
    5  = @(1.459612070389937, 114.8270160096505,  0.1575252929769615, 8.660232485948974, 125.1043858589664)
    6  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    7  = @(0.6753301551942219, 1.667794583268128, 337.1190423067083,  8.660232485948974, 348.1223995311196)
    8  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    9  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    10 = @(1.459612070389937, 10.29869402782916, 3.900430680208489,  8.660232485948974, 24.31896926437656)
    11 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    12 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 2).Value = $v[0]  # B - TB
    $ws.Cells.Item($row, 3).Value = $v[1]  # C - d2S
    $ws.Cells.Item($row, 4).Value = $v[2]  # D - K
    $ws.Cells.Item($row, 5).Value = $v[3]  # E - IP
    $ws.Cells.Item($row, 7).Value = $v[4]  # G - sum
}
